# Update the Nanning comic-con listing workbook to the newer scrape snapshot.
#
# What changed in the source data:
#   - Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both dropped
#     their oldest listing (the 2024-05-16 "剑守中国" event) because it has
#     now passed its window, so every remaining row shifts up by one
#     position. The running index in column A keeps counting 1..N from the
#     top, and the "想去人数" (interest count) in column F was re-scraped and
#     is slightly higher than it was the last time each row was fetched.
#   - Sheet "演出" (Performances) just got a refreshed interest count for its
#     single remaining row.
#   - Sheet "本地生活" (Local life) has no data rows and is untouched.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions): 14 rows -> 13 rows -------------------
$wsExpo = $wb.Worksheets.Item("展览")

# Drop the first data row (row 2); everything below slides up one row,
# carrying its formatting/number-format/shared-string type along with it.
$wsExpo.Rows.Item(2).Delete()

# The running counter in column A is positional (1, 2, 3, ...), not the
# original row's id, so after the shift it must read 1..12 again rather
# than keep the old row's shifted-up number.
for ($r = 2; $r -le 13; $r++) {
    $wsExpo.Cells.Item($r, 1).Value = $r - 1
}

# Refresh the "interest count" (column F) values that moved up a row to
# match the newly re-scraped numbers.
$wsExpo.Cells.Item(2, 6).Value = 334
$wsExpo.Cells.Item(4, 6).Value = 22
$wsExpo.Cells.Item(5, 6).Value = 3230
$wsExpo.Cells.Item(6, 6).Value = 2102
$wsExpo.Cells.Item(7, 6).Value = 402
$wsExpo.Cells.Item(8, 6).Value = 150
$wsExpo.Cells.Item(9, 6).Value = 16
$wsExpo.Cells.Item(10, 6).Value = 1193
$wsExpo.Cells.Item(11, 6).Value = 215
$wsExpo.Cells.Item(12, 6).Value = 1140
$wsExpo.Cells.Item(13, 6).Value = 93

# ---- Sheet "演出" (Performances): refresh interest count --------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value = 33

# ---- Sheet "全部类型" (All types): 15 rows -> 14 rows ------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Rows.Item(2).Delete()

for ($r = 2; $r -le 14; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$wsAll.Cells.Item(2, 6).Value = 334
$wsAll.Cells.Item(4, 6).Value = 22
$wsAll.Cells.Item(5, 6).Value = 3230
$wsAll.Cells.Item(6, 6).Value = 2102
$wsAll.Cells.Item(7, 6).Value = 402
$wsAll.Cells.Item(8, 6).Value = 33
$wsAll.Cells.Item(9, 6).Value = 150
$wsAll.Cells.Item(10, 6).Value = 16
$wsAll.Cells.Item(11, 6).Value = 1193
$wsAll.Cells.Item(12, 6).Value = 215
$wsAll.Cells.Item(13, 6).Value = 1140
$wsAll.Cells.Item(14, 6).Value = 93
